$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove column A ("Nota"), shifting B:E left to A:D
$ws.Range("A1").EntireColumn.Delete()

# Add the new descriptive row beneath the headers
$ws.Range("A2").Value = "Respondeu 99.9% das reclamações recebidas."
$ws.Range("B2").Value = "Dos que avaliaram, 62.8% voltariam a fazer negócio."
$ws.Range("C2").Value = "A empresa resolveu 67.7% das reclamações recebidas."
$ws.Range("D2").Value = "O consumidor avaliou o atendimento dessa empresa como BOM. A nota média nos últimos 6 meses é 7.0/10."
